$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) entirely, shifting rows up
$ws.Rows.Item(26).Delete()

# Delete the "SC 92" row (now row 27 after the first delete), shifting rows up again
$ws.Rows.Item(27).Delete()

# Fill in previously-missing value for "SC 5" row (now row 26), column C (B header)
$ws.Range("C26").Value = 10.8

# Clear the value for "SC 101" row (now row 27), column C (B header) - now missing
$ws.Range("C27").Value = ""

# Fill in previously-missing value for "SC 232" row (now row 33), column F
$ws.Range("F33").Value = 17.53
